# Fixed inversed condition in UltimateScalper.
# Update existing Test #1 row (row 2) and add three new UltimateScalper
# test rows (3, 4, 5) on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 2: correct the previously inverted/incorrect test case ---
$ws.Range("C2").Value = "BTCUSDT"
$ws.Range("D2").Value = 44197
$ws.Range("F2").Value = "3m"
$ws.Range("G2").Value = 1
$ws.Range("I2").Value = "UltimateScalper"

# --- Row 3: new test case (BTCUSDT, 5m) ---
$ws.Range("A3:K3").Copy() | Out-Null
$ws.Range("A3:K3").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Bybit"
$ws.Range("C3").Value = "BTCUSDT"
$ws.Range("D3").Value = 44197
$ws.Range("E3").Value = 44926
$ws.Range("F3").Value = "5m"
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 0.8
$ws.Range("I3").Value = "UltimateScalper"
$ws.Range("J3").Value = "FixedPCT"
$ws.Rows.Item(3).RowHeight = 14.25

# --- Row 4: new test case (ETHUSDT, 3m) ---
$ws.Range("A4:K4").Copy() | Out-Null
$ws.Range("A4:K4").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Bybit"
$ws.Range("C4").Value = "ETHUSDT"
$ws.Range("D4").Value = 44197
$ws.Range("E4").Value = 44926
$ws.Range("F4").Value = "3m"
$ws.Range("G4").Value = 1
$ws.Range("H4").Value = 0.8
$ws.Range("I4").Value = "UltimateScalper"
$ws.Range("J4").Value = "FixedPCT"
$ws.Rows.Item(4).RowHeight = 14.25

# --- Row 5: new test case (ETHUSDT, 5m) ---
$ws.Range("A4:K4").Copy() | Out-Null
$ws.Range("A5:K5").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "Bybit"
$ws.Range("C5").Value = "ETHUSDT"
$ws.Range("D5").Value = 44197
$ws.Range("E5").Value = 44926
$ws.Range("F5").Value = "5m"
$ws.Range("G5").Value = 1
$ws.Range("H5").Value = 0.8
$ws.Range("I5").Value = "UltimateScalper"
$ws.Range("J5").Value = "FixedPCT"
$ws.Rows.Item(5).RowHeight = 14.25

# --- Update the active selection to reflect where the user left off ---
$ws.Range("D8").Select() | Out-Null
